$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 16.4607
$ws.Range("D3").Value = -8.606299999999997
$ws.Range("E3").Value = 16.4471
$ws.Range("D4").Value = -7.520000000000005
$ws.Range("E9").Value = 17.4412
$ws.Range("B11").Value = 5.659600000000003
$ws.Range("B12").Value = 4.863399999999998
$ws.Range("D14").Value = -7.344400000000003
$ws.Range("B15").Value = 5.282299999999997
$ws.Range("E15").Value = 16.5492
$ws.Range("E19").Value = 16.3011
$ws.Range("E20").Value = 16.08379999999999
$ws.Range("E25").Value = 17.23970000000001
$ws.Range("D26").Value = -8.8422
$ws.Range("B27").Value = 6.806500000000002
$ws.Range("E27").Value = 16.55439999999999
$ws.Range("B28").Value = 6.160900000000002
$ws.Range("E28").Value = 16.81
$ws.Range("E30").Value = 15.6268
$ws.Range("B31").Value = 6.17
$ws.Range("D31").Value = -8.770399999999997
$ws.Range("B32").Value = 6.329400000000001
$ws.Range("E32").Value = 15.5584
$ws.Range("D35").Value = -9.164599999999993
$ws.Range("B36").Value = 9.795700000000005
$ws.Range("D37").Value = -7.926499999999995
$ws.Range("B38").Value = 5.270599999999999
$ws.Range("D39").Value = -7.402200000000007
$ws.Range("D40").Value = -7.957999999999995
$ws.Range("E44").Value = 16.11669999999999
$ws.Range("D45").Value = -7.604100000000003
$ws.Range("B46").Value = 7.253599999999998
$ws.Range("E47").Value = 16.5575
$ws.Range("D52").Value = -7.636099999999996
$ws.Range("B54").Value = 4.752499999999999
$ws.Range("B55").Value = 5.257599999999999
$ws.Range("B56").Value = 5.161799999999999
$ws.Range("D57").Value = -8.5403
$ws.Range("E58").Value = 16.8529
$ws.Range("E62").Value = 16.50370000000001
$ws.Range("B67").Value = 5.579099999999999
$ws.Range("B69").Value = 5.130400000000001
$ws.Range("B72").Value = 5.207500000000002
$ws.Range("B73").Value = 8.266399999999997
$ws.Range("E77").Value = 17.11690000000002
$ws.Range("E78").Value = 16.51090000000002
$ws.Range("D81").Value = -7.022099999999995
$ws.Range("B83").Value = 5.511799999999995
$ws.Range("D83").Value = -9.035199999999996
$ws.Range("E84").Value = 16.5579
$ws.Range("B86").Value = 5.117100000000002
$ws.Range("E89").Value = 17.35310000000001
$ws.Range("B91").Value = 5.415800000000001
$ws.Range("E91").Value = 17.99790000000001
$ws.Range("E92").Value = 18.09380000000002
$ws.Range("B93").Value = 5.786900000000002
$ws.Range("E96").Value = 16.0136
$ws.Range("B99").Value = 4.915299999999998
$ws.Range("D100").Value = -8.547700000000004
$ws.Range("D102").Value = -7.441499999999998
$ws.Range("E102").Value = 16.68990000000001
